$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overordnet projektplan")

# Insert a new column at Q (shifts Q:AF -> R:AG)
$ws.Columns("Q:Q").Insert()

# Fill in the new task data that was typed after the column insert.
# Row 6: "Rapport skrivning" (N6) now has hours in O6, followed by a new task
# "review af gruppe 2" (P6) with its hours "1,5 timer" (Q6).
$ws.Range("O6").Value = "5 timer"
$ws.Range("P6").Value = "review af gruppe 2"
$ws.Range("Q6").Value = "1,5 timer"

# Row 7: new task "Klargøring af fremlæggelse" (P7)
$ws.Range("P7").Value = "Klargøring af fremlæggelse"

$ws.Range("Q7").Select()
